$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for A2:G51 (Index, Time, SegStart, SegEnd, AWC, TVN, CTC)
# representing the "top 100 by adult mean word count" snippet selection.
$data = @(
        @(25, "5:37 AM", 729.045, 759.045, 11.335, 0, 0),
        @(66, "5:58 AM", 1964.125, 1994.125, 1.13, 0, 0.5),
        @(88, "6:09 AM", 2629.6425, 2659.6425, 5.215, 0, 0.25),
        @(97, "6:13 AM", 2882.4, 2912.4, 19.2, 0, 2),
        @(125, "6:27 AM", 3732.76, 3762.76, 14.64, 0, 1),
        @(139, "6:34 AM", 4149.392, 4179.392, 1.24, 0, 0.2),
        @(152, "6:41 AM", 4544.685, 4574.685, 4.463333, 0, 0.166667),
        @(161, "6:45 AM", 4808.78, 4838.78, 2.95, 0, 0),
        @(174, "6:52 AM", 5193.04, 5223.04, 71.97, 0, 0),
        @(182, "6:56 AM", 5440.89, 5470.89, 108.74, 0, 0),
        @(236, "7:23 AM", 7065.746667, 7095.746667, 0.863333, 0, 0.333333),
        @(248, "7:29 AM", 7417.6925, 7447.6925, 3.91, 0, 0),
        @(256, "7:33 AM", 7665.623333, 7695.623333, 5.7, 0.346667, 0),
        @(290, "7:50 AM", 8698.795, 8728.795, 1.96, 0, 0),
        @(303, "7:56 AM", 9072.33, 9102.33, 3.06, 0, 0.5),
        @(349, "8:19 AM", 10459.5, 10489.5, 1.52, 0, 0),
        @(362, "8:26 AM", 10841.275, 10871.275, 4.755, 0, 0.5),
        @(372, "8:31 AM", 11144.155, 11174.155, 5.175, 0, 0),
        @(381, "8:35 AM", 11413.92, 11443.92, 1.94, 19.395, 0),
        @(398, "8:44 AM", 11919.3225, 11949.3225, 1, 2.18, 0.25),
        @(515, "9:43 AM", 15448.26, 15478.26, 2.15, 0, 1),
        @(550, "10:00 A", 16484, 16514, 26.08, 0, 1),
        @(980, "1:35 PM", 29392.66, 29422.66, 6.8, 0, 0),
        @(987, "1:39 PM", 29607.37, 29637.37, 3.695, 0, 0),
        @(998, "1:44 PM", 29923.235, 29953.235, 7.655, 0, 1),
        @(1034, "2:02 PM", 31003.725, 31033.725, 8.225, 0, 0),
        @(1054, "2:12 PM", 31608.545, 31638.545, 1.995, 0, 0),
        @(1064, "2:17 PM", 31902.553333, 31932.553333, 10.7, 0, 0),
        @(1085, "2:27 PM", 32536.5, 32566.5, 0.8375, 0, 0),
        @(1092, "2:31 PM", 32747.78, 32777.78, 1.665, 0, 0.5),
        @(1103, "2:36 PM", 33070.88, 33100.88, 25.1, 0, 1),
        @(1172, "3:11 PM", 35139.775, 35169.775, 14.545, 0, 0),
        @(1195, "3:22 PM", 35827.2775, 35857.2775, 4.8075, 0, 0),
        @(1243, "3:46 PM", 37267.225, 37297.225, 9.237500000000001, 0, 0),
        @(1254, "3:52 PM", 37603.08, 37633.08, 3.615, 0, 0.25),
        @(1263, "3:57 PM", 37889.45, 37919.45, 3.78, 0, 0),
        @(1271, "4:00 PM", 38112.49, 38142.49, 8.436667, 0, 0),
        @(1297, "4:13 PM", 38885.46, 38915.46, 73.44, 0, 0),
        @(1308, "4:19 PM", 39217.913333, 39247.913333, 29.03, 0, 0),
        @(1335, "4:32 PM", 40027.72, 40057.72, 47.705, 0, 0.5),
        @(1360, "4:45 PM", 40775.5, 40805.5, 27.96, 0, 1.5),
        @(1370, "4:50 PM", 41074.77, 41104.77, 24.045, 0, 1),
        @(1393, "5:01 PM", 41774.303333, 41804.303333, 14.303333, 0, 0.333333),
        @(1401, "5:06 PM", 42025.86, 42055.86, 3.2, 0, 2),
        @(1416, "5:13 PM", 42467.672, 42497.672, 1.524, 0.2, 0),
        @(1432, "5:21 PM", 42943.82, 42973.82, 31.495, 0, 0),
        @(1445, "5:27 PM", 43323.345, 43353.345, 29.31, 0.505, 1),
        @(1482, "5:46 PM", 44444.21, 44474.21, 10.0225, 0.5024999999999999, 1.25),
        @(1534, "6:12 PM", 46012.16, 46042.16, 2.17, 0, 0),
        @(1546, "6:18 PM", 46356.505, 46386.505, 25.355, 0, 0.5)
    )

$rowCount = $data.Count
$colCount = 7

$values = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rangeAddr = "A" + $startRow + ":G" + $endRow
$range = $ws.Range($rangeAddr)
$range.Value = $values
